$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Story Status for rows that moved from "To Do" to "In Progress"
$ws.Range("E7").Value = "In Progress"
$ws.Range("E11").Value = "In Progress"
$ws.Range("E12").Value = "In Progress"
$ws.Range("E13").Value = "In Progress"
$ws.Range("E14").Value = "In Progress"
$ws.Range("E15").Value = "In Progress"
$ws.Range("E16").Value = "In Progress"
$ws.Range("E17").Value = "In Progress"
$ws.Range("E18").Value = "In Progress"
$ws.Range("E22").Value = "In Progress"
$ws.Range("E23").Value = "In Progress"

# Set Sprint Number for row 10
$ws.Range("C10").Value = 5

# Update selection to match the new view
$ws.Range("C5").Select()
